# Fruta / hortaliza, semanal
# Insert a new weekly record for "Espárragos" at row 38, pushing the
# existing rows (38-81) down to (39-82).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 5
$ws.Range("B38").Value = "Macroferia Regional de Talca"
$ws.Range("C38").Value = "Maule"
$ws.Range("D38").Value = 44874
$ws.Range("E38").Value = 7
$ws.Range("F38").Value = 300000000
$ws.Range("G38").Value = "Espárragos"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 1100
$ws.Range("L38").Value = 1100
$ws.Range("M38").Value = 1100
$ws.Range("N38").Value = '$/kilo'
$ws.Range("O38").Value = "Provincia de Linares"
$ws.Range("P38").Value = 1100
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = "Hortaliza"
